$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 0.028814
$ws.Range("N2").Value = 0.086442
$ws.Range("O2").Value = 0.003707384188741118
$ws.Range("P2").Value = 0.003707384188741118
$ws.Range("Q2").Value = 0.012166970826
$ws.Range("R2").Value = 0.109502737434
$ws.Range("S2").Value = 0.003707384188741118
$ws.Range("T2").Value = 0.003707384188741118

# Row 3 updates
$ws.Range("O3").Value = 0.9492453949340737
$ws.Range("P3").Value = 0.9492453949340736
$ws.Range("S3").Value = 0.9492453949340737
$ws.Range("T3").Value = 0.9492453949340736

# Row 4 updates
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3656536666666667
$ws.Range("N4").Value = 1.096961
$ws.Range("O4").Value = 0.04704722087718523
$ws.Range("P4").Value = 0.04704722087718523
$ws.Range("Q4").Value = 0.154400551633
$ws.Range("R4").Value = 1.389604964697
$ws.Range("S4").Value = 0.04704722087718523
$ws.Range("T4").Value = 0.04704722087718523
